# Apply "Legs Update and Sesi 2 Update!" edits to the calibration calculator.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data updates (Legs Update) ---
# Row 5 ("Leg" entry): D5 1570 -> 1580, G5 2000 -> 2050
# Dependent formulas in H5 (=D5-2*(D5-G5)) and K5 (=D5-G5) recalc automatically.
$ws.Range("D5").Value = 1580
$ws.Range("G5").Value = 2050

# Row 8: D8 1200 -> 1250
# Dependent formulas in H8 (=D8-2*(D8-G8)) and K8 (=D8-G8) recalc automatically.
$ws.Range("D8").Value = 1250

# --- Selection change on Sheet1 (G11 -> F14) ---
$ws.Range("F14").Select()

# --- Workbook window position (xWindow 6810 -> 5295) ---
# Reflected on the window object for completeness, even though the window
# chrome position is cosmetic metadata.
$win = $excel.ActiveWindow
$win.Left = 5295
